# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on row 2 of the
# zh-cn and de-de language sheets to reflect the newly generated handback
# report timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-13 21:16:03"
$zhcn.Range("H2").Value = "2016-03-13 21:16:27"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-13 21:16:07"
$dede.Range("H2").Value = "2016-03-13 21:16:33"
